# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (same 7-column fund-holding layout as
#    the other quarterly sheets) right before the "总计" (totals) sheet.
# 2. Add a new top data row to "总计" summarising the 2022-Q1 quarter, and
#    shift/renumber the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet.
# "2021-Q4" already has the exact target layout/styling (header labels,
# bold/bordered header row + index column), so clone it and just swap in
# the new quarter's numbers. Copy($before) places the clone immediately
# before the given sheet, i.e. right before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($totalSheet)

$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Fund 160416 (row 2) - code/name unchanged, only the figures move.
$q1.Range("D2").Value = "'3.37"
$q1.Range("E2").Value = "'95.08"
$q1.Range("F2").Value = "'6.39"
$q1.Range("G2").Value = "'0.2153"
$q1.Range("H2").Value = 3

# Fund 378006 (row 3) - code/name unchanged, only the figures move.
$q1.Range("D3").Value = "'0.46"
$q1.Range("E3").Value = "'88.99"
$q1.Range("F3").Value = "'2.87"
$q1.Range("G3").Value = "'0.0132"
$q1.Range("H3").Value = 4

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet with a new leading row for 2022-Q1.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

$tot.Rows.Item(2).Insert()
# Insert() drags a border/bold style down from the header row onto the
# new blank row - strip that back to the plain data-row look.
$tot.Range("B2:D2").ClearFormats()
# ...but column A keeps the bold/bordered "index" style on every data
# row, so copy that formatting (only) from the row below onto A2.
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 2
$tot.Range("D2").Value = 0.23

# Renumber the index column for the rows that shifted down.
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5
